$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 156772
$ws.Range("C4").Value = 147849
$ws.Range("C5").Value = 8923
$ws.Range("C8").Value = 63.92
